$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents() | Out-Null
$ws.Range("H12").Value = 778
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("H17").Value = 923.4737
$ws.Range("J17").Value = 923.4737
$ws.Range("L17").Value = 2770.4211
$ws.Range("N17").Value = -3106.4211
$ws.Range("H28").Value = 1459.2222
$ws.Range("I28").Value = 1599.25
$ws.Range("K28").Value = 1599.25
$ws.Range("M28").Value = -1114.25
$ws.Range("H43").Value = 1022
$ws.Range("I43").Value = 970
$ws.Range("K43").Value = 970
$ws.Range("M43").Value = -901
$ws.Range("H53").Value = 303.5
$ws.Range("I53").Value = 222
$ws.Range("J53").Value = 385
$ws.Range("K53").Value = 222
$ws.Range("L53").Value = 385
$ws.Range("M53").Value = 415
$ws.Range("N53").Value = -1659
$ws.Range("H62").Value = 1646
$ws.Range("I62").Value = 405
$ws.Range("J62").Value = 2887
$ws.Range("K62").Value = 405
$ws.Range("L62").Value = 2887
$ws.Range("M62").Value = 219
$ws.Range("N62").Value = -4135
$ws.Range("H64").Value = 4833.3335
$ws.Range("H65").Value = 1646
$ws.Range("I65").Value = 405
$ws.Range("J65").Value = 2887
$ws.Range("K65").Value = 2025
$ws.Range("L65").Value = 14435
$ws.Range("M65").Value = 1095
$ws.Range("N65").Value = -20675
$ws.Range("H67").Value = 4833.3335
$ws.Range("H86").Value = 145005.42
$ws.Range("I86").Value = 2505.6667
$ws.Range("K86").Value = 2505.6667
$ws.Range("M86").Value = -1382.6667
$ws.Range("H89").Value = 145005.42
$ws.Range("I89").Value = 2505.6667
$ws.Range("K89").Value = 12528.3335
$ws.Range("M89").Value = -6912.333500000001
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents() | Out-Null
$ws.Range("H116").Value = 10450.333
$ws.Range("I116").Value = 17952.5
$ws.Range("K116").Value = 17952.5
$ws.Range("M116").Value = -14510.5
$ws.Range("H125").Value = 10949.714
$ws.Range("I125").Value = 1412
$ws.Range("K125").Value = 12708
$ws.Range("M125").Value = -10248
$ws.Range("H137").Value = 1332.4375
$ws.Range("I137").Value = 1278.0834
$ws.Range("K137").Value = 3834.2502
$ws.Range("M137").Value = -1284.2502

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5000
$ws.Range("I110").Value = 5000
$ws.Range("K110").Value = 5000
$ws.Range("M110").Value = -2955
$ws.Range("H130").Value = 52494.5
$ws.Range("J130").Value = 52494.5
$ws.Range("L130").Value = 52494.5
$ws.Range("N130").Value = -62534.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2350.1667
$ws.Range("I86").Value = 1820.2
$ws.Range("K86").Value = 1820.2
$ws.Range("M86").Value = -697.2
$ws.Range("H89").Value = 2350.1667
$ws.Range("I89").Value = 1820.2
$ws.Range("K89").Value = 9101
$ws.Range("M89").Value = -3485

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3638.2
$ws.Range("I16").Value = 1945.1
$ws.Range("J16").Value = 7024.4
$ws.Range("K16").Value = 1945.1
$ws.Range("L16").Value = 7024.4
$ws.Range("M16").Value = -1658.1
$ws.Range("N16").Value = -7598.4
$ws.Range("H19").Value = 115.875
$ws.Range("I19").Value = 118.14286
$ws.Range("K19").Value = 118.14286
$ws.Range("M19").Value = 51.85714
$ws.Range("H24").Value = 115.875
$ws.Range("I24").Value = 118.14286
$ws.Range("K24").Value = 118.14286
$ws.Range("M24").Value = 51.85714
$ws.Range("H31").Value = 1915.375
$ws.Range("I31").Value = 1915.375
$ws.Range("K31").Value = 1915.375
$ws.Range("M31").Value = -1620.375
$ws.Range("H34").Value = 1915.375
$ws.Range("I34").Value = 1915.375
$ws.Range("K34").Value = 1915.375
$ws.Range("M34").Value = -1713.375
$ws.Range("H62").Value = 3799
$ws.Range("I62").Value = 3799
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3799
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3175
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 3799
$ws.Range("I65").Value = 3799
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 18995
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -15875
$ws.Range("N65").ClearContents() | Out-Null
$ws.Range("H86").Value = 12049.444
$ws.Range("I86").Value = 12409.833
$ws.Range("K86").Value = 12409.833
$ws.Range("M86").Value = -11286.833
$ws.Range("H89").Value = 12049.444
$ws.Range("I89").Value = 12409.833
$ws.Range("K89").Value = 62049.165
$ws.Range("M89").Value = -56433.165
$ws.Range("H109").Value = 41950
$ws.Range("J109").Value = 41950
$ws.Range("L109").Value = 41950
$ws.Range("N109").Value = -44030
$ws.Range("H113").Value = 3638.2
$ws.Range("I113").Value = 1945.1
$ws.Range("J113").Value = 7024.4
$ws.Range("K113").Value = 1945.1
$ws.Range("L113").Value = 7024.4
$ws.Range("M113").Value = 224.9000000000001
$ws.Range("N113").Value = -11364.4
$ws.Range("H132").Value = 4175.8
$ws.Range("J132").Value = 2965.3333
$ws.Range("L132").Value = 8895.999899999999
$ws.Range("N132").Value = -13955.9999
$ws.Range("H134").Value = 2341.2
$ws.Range("I134").Value = 2349.7144
$ws.Range("K134").Value = 7049.1432
$ws.Range("M134").Value = -4514.1432
$ws.Range("H141").Value = 36666.125
$ws.Range("J141").Value = 36666.125
$ws.Range("L141").Value = 36666.125
$ws.Range("N141").Value = -47026.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.05556
$ws.Range("I2").Value = 25.444445
$ws.Range("K2").Value = 152.66667
$ws.Range("M2").Value = -39.66667000000001
$ws.Range("H4").Value = 6875121.5
$ws.Range("I4").Value = 7857256
$ws.Range("J4").Value = 180
$ws.Range("K4").Value = 23571768
$ws.Range("L4").Value = 540
$ws.Range("M4").Value = -23571656
$ws.Range("N4").Value = -764
$ws.Range("H11").Value = 154
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 190
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 570
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -850
$ws.Range("H12").Value = 170.72728
$ws.Range("J12").Value = 132
$ws.Range("L12").Value = 396
$ws.Range("N12").Value = -742
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 100
$ws.Range("K26").Value = 300
$ws.Range("M26").Value = -12
$ws.Range("H64").Value = 745
$ws.Range("J64").Value = 745
$ws.Range("L64").Value = 2235
$ws.Range("N64").Value = -2775
$ws.Range("H67").Value = 745
$ws.Range("J67").Value = 745
$ws.Range("L67").Value = 2235
$ws.Range("N67").Value = -4107
$ws.Range("H137").Value = 3729.7778
$ws.Range("I137").Value = 2395
$ws.Range("J137").Value = 3896.625
$ws.Range("K137").Value = 7185
$ws.Range("L137").Value = 11689.875
$ws.Range("M137").Value = -2085
$ws.Range("N137").Value = -21889.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19990
$ws.Range("J58").Value = 19990
$ws.Range("L58").Value = 19990
$ws.Range("N58").Value = -20544
$ws.Range("H113").Value = 798
$ws.Range("I113").Value = 497.2
$ws.Range("K113").Value = 497.2
$ws.Range("M113").Value = 1672.8
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents() | Out-Null

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5387.25
$ws.Range("I40").Value = 5119.6
$ws.Range("K40").Value = 5119.6
$ws.Range("M40").Value = -4983.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 191869.4
$ws.Range("J2").Value = 191869.4
$ws.Range("L2").Value = 191869.4
$ws.Range("N2").Value = -192093.4
$ws.Range("H14").Value = 6514.143
$ws.Range("I14").Value = 900.6667
$ws.Range("J14").Value = 10724.25
$ws.Range("K14").Value = 900.6667
$ws.Range("L14").Value = 10724.25
$ws.Range("M14").Value = -732.6667
$ws.Range("N14").Value = -11060.25
$ws.Range("H30").Value = 3999.5
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents() | Out-Null
$ws.Range("H94").Value = 90000
$ws.Range("J94").Value = 90000
$ws.Range("L94").Value = 90000
$ws.Range("N94").Value = -91802
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 14500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 14500
$ws.Range("M95").ClearContents() | Out-Null
$ws.Range("N95").Value = -19992
$ws.Range("H96").Value = 2650.75
$ws.Range("J96").Value = 2533.3333
$ws.Range("L96").Value = 2533.3333
$ws.Range("N96").Value = -5279.3333
$ws.Range("H100").Value = 12500418
$ws.Range("I100").Value = 14286092
$ws.Range("K100").Value = 28572184
$ws.Range("M100").Value = -28571643
$ws.Range("H122").Value = 1827.7142
$ws.Range("I122").Value = 1749.1666
$ws.Range("J122").Value = 2299
$ws.Range("K122").Value = 5247.4998
$ws.Range("L122").Value = 6897
$ws.Range("M122").Value = -2797.4998
$ws.Range("N122").Value = -11797
